# Update cryptos list price/volume columns with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.626.47'
$ws.Range("D3").Value = '1.675.80'
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("E4").Value = '  -0.34%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '219.96'
$cell.ClearFormats()
$ws.Range("E5").Value = '  +2.65%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.529'
$cell.ClearFormats()
$ws.Range("E6").Value = '  +2.26%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.ClearFormats()
$ws.Range("E7").Value = '  -0.35%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '29.64'
$cell.ClearFormats()
$ws.Range("E8").Value = '  +4.84%  '
$ws.Range("E9").Value = '  +2.99%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0648'
$cell.ClearFormats()
$ws.Range("E10").Value = '  +6.69%  '
$ws.Range("D12").Value = '1.913.94'
$ws.Range("E12").Value = '  +2.66%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '10.29'
$cell.ClearFormats()
$ws.Range("E13").Value = '  +13.08%  '
$ws.Range("D14").Value = '1.674.01'
$ws.Range("E14").Value = '  +2.71%  '
$ws.Range("E15").Value = '  +9.46%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '4.02'
$cell.ClearFormats()
$ws.Range("E16").Value = '  +4.78%  '
$ws.Range("D17").Value = '30.632.06'
$ws.Range("E17").Value = '  +2.33%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '66.45'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +3.94%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '244.09'
$cell.ClearFormats()
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = '0.0₃0725'
$ws.Range("E20").Value = '  +3.44%  '
$ws.Range("E21").Value = '  -0.28%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '10.02'
$cell.ClearFormats()
$ws.Range("E23").Value = '  +2.98%  '
$ws.Range("E24").Value = '  +0.40%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '158.62'
$cell.ClearFormats()
$ws.Range("E25").Value = '  -0.91%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '15.90'
$cell.ClearFormats()
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("E27").Value = '  +2.88%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '6.70'
$cell.ClearFormats()
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("E30").Value = '  +2.29%  '
$ws.Range("E31").Value = '  +3.82%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.47'
$cell.ClearFormats()
$ws.Range("E32").Value = '  +3.26%  '
$ws.Range("E33").Value = '  +4.38%  '
$ws.Range("D34").Value = '1.482.87'
$ws.Range("E34").Value = '  +4.15%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.77'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +7.73%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '84.75'
$cell.ClearFormats()
$ws.Range("E36").Value = '  +12.59%  '
$ws.Range("E37").Value = '  -0.41%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.602'
$cell.ClearFormats()
$ws.Range("E38").Value = '  +9.51%  '
$ws.Range("E39").Value = '  +5.73%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.67'
$cell.ClearFormats()
$ws.Range("E40").Value = '  -3.26%  '
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("E43").Value = '  +1.64%  '
$ws.Range("E44").Value = '  -0.51%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell.ClearFormats()
$ws.Range("E45").Value = '  -0.13%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -0.33%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '51.45'
$cell.ClearFormats()
$ws.Range("E48").Value = '  +3.29%  '
$ws.Range("E49").Value = '  +1.99%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '94.84'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +5.02%  '
$ws.Range("E51").Value = '  +0.02%  '
